# TC46_Canine_Filter_Breed-WestHlnd.xlsx - "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab query (row 4, column B on the "startup" sheet) dropped the
# `File Type` and `Breed` columns from its Neo4j RETURN clause (Diagnosis
# keeps one extra leading space where the removed Breed line used to be).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['West Highland White Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# The query text got shorter, so the wrapped row shrinks from 246.5 to 217.5.
$ws.Rows.Item(4).RowHeight = 217.5

# Scroll/select so row 4 is in view, matching the saved cursor position.
$null = $ws.Range("B4").Select()
